# Applies the weekly Fruta/Hortaliza update for "Espárragos" (Mercado Mayorista
# Lo Valledor de Santiago) rows 93-112:
#  - rows 93-108: updated figures for existing weekly records (dates shifted
#    down one week, new 2021-11-23 week inserted, prices/quality/origin refreshed)
#  - rows 109-112: four new records appended (two carried forward from the old
#    week, two reinstating the original "Verde" variety rows that moved down)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 93
$ws.Cells.Item(93, 4).Value = 44523
$ws.Cells.Item(93, 10).Value = 410
$ws.Cells.Item(93, 11).Value = 1300
$ws.Cells.Item(93, 12).Value = 1400
$ws.Cells.Item(93, 13).Value = 1356
$ws.Cells.Item(93, 15).Value = "Región Metropolitana"
$ws.Cells.Item(93, 16).Value = 1356

# Row 94
$ws.Cells.Item(94, 4).Value = 44523
$ws.Cells.Item(94, 10).Value = 330
$ws.Cells.Item(94, 11).Value = 1100
$ws.Cells.Item(94, 12).Value = 1200
$ws.Cells.Item(94, 13).Value = 1155
$ws.Cells.Item(94, 15).Value = "Región Metropolitana"
$ws.Cells.Item(94, 16).Value = 1155

# Row 95
$ws.Cells.Item(95, 4).Value = 44523
$ws.Cells.Item(95, 10).Value = 220
$ws.Cells.Item(95, 11).Value = 900
$ws.Cells.Item(95, 12).Value = 1000
$ws.Cells.Item(95, 13).Value = 955
$ws.Cells.Item(95, 15).Value = "Región Metropolitana"
$ws.Cells.Item(95, 16).Value = 955

# Row 96
$ws.Cells.Item(96, 4).Value = 44474
$ws.Cells.Item(96, 10).Value = 780
$ws.Cells.Item(96, 11).Value = 1500
$ws.Cells.Item(96, 12).Value = 1600
$ws.Cells.Item(96, 13).Value = 1558
$ws.Cells.Item(96, 16).Value = 1558

# Row 97
$ws.Cells.Item(97, 4).Value = 44474
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 520
$ws.Cells.Item(97, 13).Value = 1348
$ws.Cells.Item(97, 15).Value = "Provincia de Linares"
$ws.Cells.Item(97, 16).Value = 1348

# Row 98
$ws.Cells.Item(98, 4).Value = 44474
$ws.Cells.Item(98, 9).Value = "Segunda"
$ws.Cells.Item(98, 10).Value = 400
$ws.Cells.Item(98, 11).Value = 1000
$ws.Cells.Item(98, 12).Value = 1200
$ws.Cells.Item(98, 13).Value = 1100
$ws.Cells.Item(98, 16).Value = 1100

# Row 99
$ws.Cells.Item(99, 9).Value = "Banquete"
$ws.Cells.Item(99, 10).Value = 990
$ws.Cells.Item(99, 11).Value = 1400
$ws.Cells.Item(99, 12).Value = 1500
$ws.Cells.Item(99, 13).Value = 1443
$ws.Cells.Item(99, 15).Value = "Provincia de Linares"
$ws.Cells.Item(99, 16).Value = 1443

# Row 100
$ws.Cells.Item(100, 9).Value = "Banquete"
$ws.Cells.Item(100, 10).Value = 660
$ws.Cells.Item(100, 11).Value = 1300
$ws.Cells.Item(100, 12).Value = 1400
$ws.Cells.Item(100, 13).Value = 1335
$ws.Cells.Item(100, 15).Value = "Región Metropolitana"
$ws.Cells.Item(100, 16).Value = 1335

# Row 101
$ws.Cells.Item(101, 9).Value = "Primera"
$ws.Cells.Item(101, 10).Value = 1000
$ws.Cells.Item(101, 11).Value = 1200
$ws.Cells.Item(101, 12).Value = 1300
$ws.Cells.Item(101, 13).Value = 1253
$ws.Cells.Item(101, 15).Value = "Provincia de Linares"
$ws.Cells.Item(101, 16).Value = 1253

# Row 102
$ws.Cells.Item(102, 4).Value = 44491
$ws.Cells.Item(102, 9).Value = "Primera"
$ws.Cells.Item(102, 10).Value = 640
$ws.Cells.Item(102, 11).Value = 1100
$ws.Cells.Item(102, 12).Value = 1200
$ws.Cells.Item(102, 13).Value = 1142
$ws.Cells.Item(102, 15).Value = "Región Metropolitana"
$ws.Cells.Item(102, 16).Value = 1142

# Row 103
$ws.Cells.Item(103, 4).Value = 44491
$ws.Cells.Item(103, 9).Value = "Segunda"
$ws.Cells.Item(103, 10).Value = 630
$ws.Cells.Item(103, 11).Value = 1000
$ws.Cells.Item(103, 12).Value = 1000
$ws.Cells.Item(103, 13).Value = 1000
$ws.Cells.Item(103, 15).Value = "Provincia de Linares"
$ws.Cells.Item(103, 16).Value = 1000

# Row 104
$ws.Cells.Item(104, 4).Value = 44491
$ws.Cells.Item(104, 9).Value = "Segunda"
$ws.Cells.Item(104, 10).Value = 440
$ws.Cells.Item(104, 11).Value = 900
$ws.Cells.Item(104, 12).Value = 1000
$ws.Cells.Item(104, 13).Value = 941
$ws.Cells.Item(104, 15).Value = "Región Metropolitana"
$ws.Cells.Item(104, 16).Value = 941

# Row 105
$ws.Cells.Item(105, 9).Value = "Banquete"
$ws.Cells.Item(105, 10).Value = 1160
$ws.Cells.Item(105, 11).Value = 1300
$ws.Cells.Item(105, 12).Value = 1400
$ws.Cells.Item(105, 13).Value = 1347
$ws.Cells.Item(105, 15).Value = "Provincia de Linares"
$ws.Cells.Item(105, 16).Value = 1347

# Row 106
$ws.Cells.Item(106, 9).Value = "Banquete"
$ws.Cells.Item(106, 10).Value = 1000
$ws.Cells.Item(106, 11).Value = 1200
$ws.Cells.Item(106, 12).Value = 1300
$ws.Cells.Item(106, 13).Value = 1268
$ws.Cells.Item(106, 15).Value = "Región Metropolitana"
$ws.Cells.Item(106, 16).Value = 1268

# Row 107
$ws.Cells.Item(107, 9).Value = "Primera"
$ws.Cells.Item(107, 10).Value = 930
$ws.Cells.Item(107, 11).Value = 1100
$ws.Cells.Item(107, 12).Value = 1200
$ws.Cells.Item(107, 13).Value = 1152
$ws.Cells.Item(107, 15).Value = "Provincia de Linares"
$ws.Cells.Item(107, 16).Value = 1152

# Row 108
$ws.Cells.Item(108, 4).Value = 44495
$ws.Cells.Item(108, 8).Value = "Sin especificar"
$ws.Cells.Item(108, 10).Value = 710
$ws.Cells.Item(108, 12).Value = 1100
$ws.Cells.Item(108, 13).Value = 1063
$ws.Cells.Item(108, 15).Value = "Región Metropolitana"
$ws.Cells.Item(108, 16).Value = 1063

# Row 109 (new)
$ws.Cells.Item(109, 1).Value = 6
$ws.Cells.Item(109, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(109, 3).Value = "Metropolitana"
$ws.Cells.Item(109, 4).Value = 44495
$ws.Cells.Item(109, 5).Value = 13
$ws.Cells.Item(109, 6).Value = 300000000
$ws.Cells.Item(109, 7).Value = "Espárragos"
$ws.Cells.Item(109, 8).Value = "Sin especificar"
$ws.Cells.Item(109, 9).Value = "Segunda"
$ws.Cells.Item(109, 10).Value = 660
$ws.Cells.Item(109, 11).Value = 900
$ws.Cells.Item(109, 12).Value = 1000
$ws.Cells.Item(109, 13).Value = 955
$ws.Cells.Item(109, 14).Value = "$/kilo"
$ws.Cells.Item(109, 15).Value = "Provincia de Linares"
$ws.Cells.Item(109, 16).Value = 955
$ws.Cells.Item(109, 17).Value = 1
$ws.Cells.Item(109, 18).Value = "Hortaliza"
$ws.Cells.Item(109, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 110 (new)
$ws.Cells.Item(110, 1).Value = 6
$ws.Cells.Item(110, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(110, 3).Value = "Metropolitana"
$ws.Cells.Item(110, 4).Value = 44495
$ws.Cells.Item(110, 5).Value = 13
$ws.Cells.Item(110, 6).Value = 300000000
$ws.Cells.Item(110, 7).Value = "Espárragos"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Segunda"
$ws.Cells.Item(110, 10).Value = 490
$ws.Cells.Item(110, 11).Value = 800
$ws.Cells.Item(110, 12).Value = 900
$ws.Cells.Item(110, 13).Value = 876
$ws.Cells.Item(110, 14).Value = "$/kilo"
$ws.Cells.Item(110, 15).Value = "Región Metropolitana"
$ws.Cells.Item(110, 16).Value = 876
$ws.Cells.Item(110, 17).Value = 1
$ws.Cells.Item(110, 18).Value = "Hortaliza"
$ws.Cells.Item(110, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 111 (new)
$ws.Cells.Item(111, 1).Value = 6
$ws.Cells.Item(111, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(111, 3).Value = "Metropolitana"
$ws.Cells.Item(111, 4).Value = 44161
$ws.Cells.Item(111, 5).Value = 13
$ws.Cells.Item(111, 6).Value = 300000000
$ws.Cells.Item(111, 7).Value = "Espárragos"
$ws.Cells.Item(111, 8).Value = "Verde"
$ws.Cells.Item(111, 9).Value = "Primera"
$ws.Cells.Item(111, 10).Value = 4300
$ws.Cells.Item(111, 11).Value = 1000
$ws.Cells.Item(111, 12).Value = 1000
$ws.Cells.Item(111, 13).Value = 1000
$ws.Cells.Item(111, 14).Value = "$/kilo"
$ws.Cells.Item(111, 15).Value = "Provincia de Linares"
$ws.Cells.Item(111, 16).Value = 1000
$ws.Cells.Item(111, 17).Value = 1
$ws.Cells.Item(111, 18).Value = "Hortaliza"
$ws.Cells.Item(111, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 112 (new)
$ws.Cells.Item(112, 1).Value = 6
$ws.Cells.Item(112, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(112, 3).Value = "Metropolitana"
$ws.Cells.Item(112, 4).Value = 44161
$ws.Cells.Item(112, 5).Value = 13
$ws.Cells.Item(112, 6).Value = 300000000
$ws.Cells.Item(112, 7).Value = "Espárragos"
$ws.Cells.Item(112, 8).Value = "Verde"
$ws.Cells.Item(112, 9).Value = "Segunda"
$ws.Cells.Item(112, 10).Value = 2500
$ws.Cells.Item(112, 11).Value = 800
$ws.Cells.Item(112, 12).Value = 800
$ws.Cells.Item(112, 13).Value = 800
$ws.Cells.Item(112, 14).Value = "$/kilo"
$ws.Cells.Item(112, 15).Value = "Provincia de Linares"
$ws.Cells.Item(112, 16).Value = 800
$ws.Cells.Item(112, 17).Value = 1
$ws.Cells.Item(112, 18).Value = "Hortaliza"
$ws.Cells.Item(112, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

